$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.914.47"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.630.49"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'215.57"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.2568"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").Value = "'0.06338"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "'19.50"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").Value = "'0.07748"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.240"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.630.81"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "1.856.20"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "'0.5482"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "0.0₅7653"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "25.905.45"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'4.423"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'194.38"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "'9.887"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "'6.049"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").Value = "'1.915"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").Value = "'141.93"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "'0.1235"
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("D28").Value = "'6.791"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "'1.239"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "'0.04879"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("D32").Value = "'3.240"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").Value = "'3.183"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("D34").Value = "'1.544"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").Value = "'2.368"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'0.8964"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.5517"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.543"
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").Value = "1.117.63"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").Value = "'0.01553"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'5.570"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'0.7992"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").Value = "'97.09"
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.768.98"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈118"
$ws.Range("E46").Value = "  -7.47%  "
$ws.Range("D47").Value = "'0.4445"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").Value = "'1.003"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").Value = "'54.76"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "'0.05142"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("D51").Value = "'7.505"
$ws.Range("E51").Value = "  +2.32%  "
